$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 208 (this shifts rows 208:280 down to 209:281)
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new data
$ws.Cells.Item(208, 1).Value = 4
$ws.Cells.Item(208, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value = "Los Lagos"
$ws.Cells.Item(208, 4).Value = 44988
$ws.Cells.Item(208, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 5).Value = 10
$ws.Cells.Item(208, 6).Value = 100112009
$ws.Cells.Item(208, 7).Value = "Acelga"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 70
$ws.Cells.Item(208, 11).Value = 10000
$ws.Cells.Item(208, 12).Value = 10000
$ws.Cells.Item(208, 13).Value = 10000
$ws.Cells.Item(208, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(208, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(208, 16).Value = 833
$ws.Cells.Item(208, 17).Value = 12
$ws.Cells.Item(208, 18).Value = "Hortaliza"
